$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.040.50"
$ws.Range("E2").Value = "  +6.70%  "
$ws.Range("D3").Value = "3.015.17"
$ws.Range("E3").Value = "  +4.06%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "583.87"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("E6").Value = "  +13.72%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.010.55"
$ws.Range("E8").Value = "  +3.97%  "
$ws.Range("E9").Value = "  +3.34%  "
$ws.Range("D10").Value = "6.98"
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").Value = "0.158"
$ws.Range("E11").Value = "  +7.85%  "
$ws.Range("E12").Value = "  +6.52%  "
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  +9.08%  "
$ws.Range("D14").Value = "34.99"
$ws.Range("E14").Value = "  +8.08%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "66.034.44"
$ws.Range("E16").Value = "  +6.77%  "
$ws.Range("D17").Value = "3.515.69"
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("E18").Value = "  +7.08%  "
$ws.Range("D19").Value = "3.013.33"
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("D20").Value = "458.57"
$ws.Range("E20").Value = "  +6.31%  "
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  +8.20%  "
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  +5.56%  "
$ws.Range("D23").Value = "7.40"
$ws.Range("E23").Value = "  +7.69%  "
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("E25").Value = "  +14.14%  "
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").Value = "10.65"
$ws.Range("E27").Value = "  +5.61%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "8.17"
$ws.Range("E29").Value = "  +16.59%  "
$ws.Range("E30").Value = "  +14.97%  "
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("E32").Value = "  -7.12%  "
$ws.Range("E33").Value = "  +5.84%  "
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +4.52%  "
$ws.Range("E37").Value = "  +7.65%  "
$ws.Range("D38").Value = "2.14"
$ws.Range("E38").Value = "  +11.73%  "
$ws.Range("D39").Value = "3.06"
$ws.Range("E39").Value = "  +7.29%  "
$ws.Range("D40").Value = "50.02"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Value = "0.309"
$ws.Range("E41").Value = "  +15.16%  "
$ws.Range("E42").Value = "  +6.40%  "
$ws.Range("D43").Value = "43.78"
$ws.Range("E43").Value = "  +8.78%  "
$ws.Range("D44").Value = "8.50"
$ws.Range("E44").Value = "  +4.43%  "
$ws.Range("D45").Value = "388.19"
$ws.Range("E45").Value = "  +11.61%  "
$ws.Range("E46").Value = "  +6.77%  "
$ws.Range("D47").Value = "2.798.59"
$ws.Range("E47").Value = "  +3.74%  "
$ws.Range("D48").Value = "135.35"
$ws.Range("E48").Value = "  +2.73%  "
$ws.Range("D50").Value = "24.08"
$ws.Range("E50").Value = "  +11.38%  "
$ws.Range("E51").Value = "  +4.17%  "
